$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (Through 2022-03-21 -> Through 2022-03-22)
$ws.Name = "Through 2022-03-22"

# Update header label in I1 (2022 (through 03-21) -> 2022 (through 03-22))
$ws.Range("I1").Value = "2022 (through 03-22)"

# Update data values for the "Total" column (I) for January, February, March, and yearly Total
$ws.Range("I2").Value = 160
$ws.Range("I3").Value = 140
$ws.Range("I4").Value = 94
$ws.Range("I14").Value = 394
